# Refresh the "ランサーズ" listing sheet with the newest scrape snapshot
# (2025-10-28 06:28:52 JST): rewrite rows 2-15 with the refreshed job
# postings/ranking, drop the now out-of-window rows 16-21, and repoint
# the F-column hyperlinks at their new URLs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Overwrite rows 2-15 with the refreshed listing data ---
$ws.Cells.Item(2, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(2, 2).Value = "医療機関向けAIアプリとLINEの連携開発を支援してくださるAIエンジニア募集(AI/バックエンド)"
$ws.Cells.Item(2, 3).Value = "システム開発"
$ws.Cells.Item(2, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(2, 5).Value = "期限情報なし"
$ws.Cells.Item(2, 6).Value = "https://www.lancers.jp/work/detail/5416301"
$ws.Cells.Item(2, 7).Value = 385
$ws.Cells.Item(2, 8).Value = "🔥AI,Ai ◆開発 ◇アプリ"

$ws.Cells.Item(3, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(3, 2).Value = "大企業の業務効率化AIプロジェクトの技術方針策定を支援するAIテックリード募集"
$ws.Cells.Item(3, 3).Value = "システム開発"
$ws.Cells.Item(3, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(3, 5).Value = "期限情報なし"
$ws.Cells.Item(3, 6).Value = "https://www.lancers.jp/work/detail/5416307"
$ws.Cells.Item(3, 7).Value = 378
$ws.Cells.Item(3, 8).Value = "🔥AI,Ai ◆効率化"

$ws.Cells.Item(4, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(4, 2).Value = "Azureでの社内文書検索RAG開発の精度改善を伴走支援してくださるAIエンジニア募集"
$ws.Cells.Item(4, 3).Value = "システム開発"
$ws.Cells.Item(4, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(4, 5).Value = "期限情報なし"
$ws.Cells.Item(4, 6).Value = "https://www.lancers.jp/work/detail/5416305"
$ws.Cells.Item(4, 7).Value = 375
$ws.Cells.Item(4, 8).Value = "🔥AI,Ai ◆開発"

$ws.Cells.Item(5, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(5, 2).Value = "Google AI studio が生成したウェブアプリの調整【AI文章の提案は受け付けません】"
$ws.Cells.Item(5, 3).Value = "システム開発"
$ws.Cells.Item(5, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(5, 5).Value = "期限情報なし"
$ws.Cells.Item(5, 6).Value = "https://www.lancers.jp/work/detail/5421873"
$ws.Cells.Item(5, 7).Value = 330
$ws.Cells.Item(5, 8).Value = "🔥AI,Ai ◇アプリ"

$ws.Cells.Item(6, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(6, 2).Value = "Stable Diffusionに詳しいLoRAなどを用いた画像生成AIエンジニア募集"
$ws.Cells.Item(6, 3).Value = "システム開発"
$ws.Cells.Item(6, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(6, 5).Value = "期限情報なし"
$ws.Cells.Item(6, 6).Value = "https://www.lancers.jp/work/detail/5416328"
$ws.Cells.Item(6, 7).Value = 310
$ws.Cells.Item(6, 8).Value = "🔥AI,Ai"

$ws.Cells.Item(7, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(7, 2).Value = "【Unity/XRエンジニア募集】製造業DX支援!既存システムと連携するXRアプリ開発"
$ws.Cells.Item(7, 3).Value = "システム開発"
$ws.Cells.Item(7, 4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(7, 5).Value = "期限情報なし"
$ws.Cells.Item(7, 6).Value = "https://www.lancers.jp/work/detail/5421687"
$ws.Cells.Item(7, 7).Value = 115
$ws.Cells.Item(7, 8).Value = "◆開発 ◇アプリ"

$ws.Cells.Item(8, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(8, 2).Value = "新卒向け就活マッチングWebサービス開発(診断・管理画面含む)"
$ws.Cells.Item(8, 3).Value = "システム開発"
$ws.Cells.Item(8, 4).Value = "5,000,000 円 ~ / 固定"
$ws.Cells.Item(8, 5).Value = "期限情報なし"
$ws.Cells.Item(8, 6).Value = "https://www.lancers.jp/work/detail/5421820"
$ws.Cells.Item(8, 7).Value = 100
$ws.Cells.Item(8, 8).Value = "◆開発 ◇管理"

$ws.Cells.Item(9, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(9, 2).Value = "【急募】Webアプリ開発エンジニア募集!フルリモート可"
$ws.Cells.Item(9, 3).Value = "システム開発"
$ws.Cells.Item(9, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(9, 5).Value = "期限情報なし"
$ws.Cells.Item(9, 6).Value = "https://www.lancers.jp/work/detail/5411585"
$ws.Cells.Item(9, 7).Value = 93
$ws.Cells.Item(9, 8).Value = "◆開発 ◇アプリ"

$ws.Cells.Item(10, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(10, 2).Value = "【急募】UberEats案件オファー抽出アプリのバックエンド開発"
$ws.Cells.Item(10, 3).Value = "システム開発"
$ws.Cells.Item(10, 4).Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Cells.Item(10, 5).Value = "期限情報なし"
$ws.Cells.Item(10, 6).Value = "https://www.lancers.jp/work/detail/5421779"
$ws.Cells.Item(10, 7).Value = 93
$ws.Cells.Item(10, 8).Value = "◆開発 ◇アプリ"

$ws.Cells.Item(11, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(11, 2).Value = "弥生販売 得意先台帳登録 商品登録 売上伝票作成ツールのご相談"
$ws.Cells.Item(11, 3).Value = "システム開発"
$ws.Cells.Item(11, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(11, 5).Value = "期限情報なし"
$ws.Cells.Item(11, 6).Value = "https://www.lancers.jp/work/detail/5422004"
$ws.Cells.Item(11, 7).Value = 73
$ws.Cells.Item(11, 8).Value = "◆ツール"

$ws.Cells.Item(12, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(12, 2).Value = "【急募】神奈川県の既存小規模オフィス向けファイル共有システム構築"
$ws.Cells.Item(12, 3).Value = "システム開発"
$ws.Cells.Item(12, 4).Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Cells.Item(12, 5).Value = "期限情報なし"
$ws.Cells.Item(12, 6).Value = "https://www.lancers.jp/work/detail/5422125"
$ws.Cells.Item(12, 7).Value = 33
$ws.Cells.Item(12, 8).ClearContents()

$ws.Cells.Item(13, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(13, 2).Value = "【システム構築】Square運用とスプレッドシート作成の支援"
$ws.Cells.Item(13, 3).Value = "システム開発"
$ws.Cells.Item(13, 4).Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Cells.Item(13, 5).Value = "期限情報なし"
$ws.Cells.Item(13, 6).Value = "https://www.lancers.jp/work/detail/5421982"
$ws.Cells.Item(13, 7).Value = 28
$ws.Cells.Item(13, 8).ClearContents()

$ws.Cells.Item(14, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(14, 2).Value = "Stable Diffusion LoRA制作依頼 画風+キャラ・衣装LoRA量産テンプレ構築"
$ws.Cells.Item(14, 3).Value = "システム開発"
$ws.Cells.Item(14, 4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(14, 5).Value = "期限情報なし"
$ws.Cells.Item(14, 6).Value = "https://www.lancers.jp/work/detail/5421894"
$ws.Cells.Item(14, 7).Value = 18
$ws.Cells.Item(14, 8).ClearContents()

$ws.Cells.Item(15, 1).Value = "2025-10-28 06:28:52"
$ws.Cells.Item(15, 2).Value = "限定公開 限定公開の仕事"
$ws.Cells.Item(15, 3).Value = "システム開発"
$ws.Cells.Item(15, 4).Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Cells.Item(15, 5).Value = "期限情報なし"
$ws.Cells.Item(15, 6).Value = "https://www.lancers.jp/work/detail/5421782"
$ws.Cells.Item(15, 7).Value = 10
$ws.Cells.Item(15, 8).ClearContents()

# --- Drop rows 16-21: they fell out of the top-15 window on this run ---
$ws.Range("A16:H21").EntireRow.Delete()

# --- Rebuild the F-column hyperlinks so each one points at its row's URL ---
# (Hyperlinks collection is sheet-scoped, so clear it out fully and re-add
# in row order to land back on a clean rId1..rId14 sequence.)
$ws.Hyperlinks.Delete()
$hyperlinkUrls = @(
    "https://www.lancers.jp/work/detail/5416301",
    "https://www.lancers.jp/work/detail/5416307",
    "https://www.lancers.jp/work/detail/5416305",
    "https://www.lancers.jp/work/detail/5421873",
    "https://www.lancers.jp/work/detail/5416328",
    "https://www.lancers.jp/work/detail/5421687",
    "https://www.lancers.jp/work/detail/5421820",
    "https://www.lancers.jp/work/detail/5411585",
    "https://www.lancers.jp/work/detail/5421779",
    "https://www.lancers.jp/work/detail/5422004",
    "https://www.lancers.jp/work/detail/5422125",
    "https://www.lancers.jp/work/detail/5421982",
    "https://www.lancers.jp/work/detail/5421894",
    "https://www.lancers.jp/work/detail/5421782"
)
for ($i = 0; $i -lt $hyperlinkUrls.Length; $i++) {
    $targetCell = $ws.Cells.Item($i + 2, 6)
    $ws.Hyperlinks.Add($targetCell, $hyperlinkUrls[$i]) | Out-Null
}

# --- Shrink the sheet dimension to match the trimmed data (A1:H15) ---
$ws.Range("A1:H15").Select()
